$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting existing D:K data to E:L
$ws.Range("D1").EntireColumn.Insert()

# Copy number formats from column E (the old column D, now shifted) into the new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate new column D with the newly-reported fiscal period figures
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 9124000
$ws.Range("D9").Value = 4627000
$ws.Range("D10").Value = 4497000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 61000
$ws.Range("D15").Value = 2230000
$ws.Range("D17").Value = 8391000
$ws.Range("D18").Value = 733000
$ws.Range("D20").Value = -259000
$ws.Range("D21").Value = 2704000
$ws.Range("D22").Value = 207000
$ws.Range("D23").Value = 267000
$ws.Range("D24").Value = 102000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 165000
$ws.Range("D27").Value = 165000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 259000
$ws.Range("D33").Value = 165000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 165000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 615000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 955000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 604000
$ws.Range("D46").Value = 2174000
$ws.Range("D47").Value = 1190000
$ws.Range("D48").Value = 12210000
$ws.Range("D49").Value = 1917000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 1658000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 19149000
$ws.Range("D57").Value = 371000
$ws.Range("D58").Value = 23000
$ws.Range("D59").Value = 1322000
$ws.Range("D60").Value = 1716000
$ws.Range("D61").Value = 13760000
$ws.Range("D62").Value = 3259000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 18735000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -1091000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 414000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 165000
$ws.Range("D83").Value = 2230000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 2609000
$ws.Range("D91").Value = -12820000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -3426000
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 667000
$ws.Range("D101").Value = -16000
$ws.Range("D102").Value = -166000

# The column insert also stamped empty styled cells into otherwise-blank rows; remove them
$ws.Range("D36").Clear()
$ws.Range("D37").Clear()
$ws.Range("D78").Clear()
$ws.Range("D79").Clear()
